# Generate Report for Handoff
# Replaces the old handoff UUID / content-hash / timestamps with the new ones
# produced by the latest localization-status generation run.

$wb = $excel.ActiveWorkbook

$oldGuid = "98169608-cdfe-4201-be01-68704be37962"
$newGuid = "7d71bcde-2188-4dea-9990-360038628121"

$oldHash = "2df7aed3e08f0fc9c725210c1b94d060d02dee8d"
$newHash = "d4ae73b8332663ba5dcae6116e1de887174e2bbf"

# NOTE: the handoff-generator only refreshes the *display* text of each
# hyperlink (it re-derives the text shown in the cell from the new file
# name) - the underlying hyperlink target (Address) still points at the
# already-published blob for the previous commit, so it is left as-is.
$oldTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1166163e10cfded37ca53ffd330078e916bea5b4/e2e/" + $oldGuid + ".md"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newGuid + ".md"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $oldTarget, [Type]::Missing, [Type]::Missing, "e2e\" + $newGuid + ".md")

$wsOverview.Range("G2").Value = "2016-09-02 07:06:44"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newGuid + ".md"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $oldTarget, [Type]::Missing, [Type]::Missing, $newGuid + ".md")

$wsZhCn.Range("G2").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-02 07:06:39"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newGuid + ".md"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $oldTarget, [Type]::Missing, [Type]::Missing, $newGuid + ".md")

$wsDeDe.Range("G2").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-02 07:06:44"
